# Update tracker.xlsx from Admin app (v1.8.8)
#
# 1) Convert the "Date" column (B) on the Events sheet from inline text
#    ("YYYY-MM-DD") to real Excel date serial numbers, formatted with the
#    existing datetime number format (reuses style index 2, already present
#    in styles.xml and already used elsewhere in the workbook).
# 2) Add a new "HighHand_Info" sheet at the end of the workbook with a
#    header row (sharing the same bold/centered/bordered header style used
#    by the other sheets) and one data row that records the last-updated
#    timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Events sheet: replace text dates in column B with real date serials
# ---------------------------------------------------------------------
$events = $wb.Worksheets.Item("Events")

$dateSerials = @{
  "B2"  = 45877   # 2025-08-08
  "B3"  = 45919   # 2025-09-19
  "B4"  = 45940   # 2025-10-10
  "B5"  = 45975   # 2025-11-14
  "B6"  = 46003   # 2025-12-12
  "B7"  = 46031   # 2026-01-09
  "B8"  = 46052   # 2026-01-30
  "B9"  = 46073   # 2026-02-20
  "B10" = 46094   # 2026-03-13
  "B11" = 46115   # 2026-04-03
  "B12" = 46136   # 2026-04-24
  "B13" = 46150   # 2026-05-08
}

foreach ($addr in $dateSerials.Keys) {
    $cell = $events.Range($addr)
    $cell.Value = $dateSerials[$addr]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# ---------------------------------------------------------------------
# 2) Add the new HighHand_Info sheet after the last existing sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$hh = $wb.Worksheets.Add($null, $lastSheet)
$hh.Name = "HighHand_Info"

# Header row text
$hh.Range("A1").Value = "Current Holder"
$hh.Range("B1").Value = "Hand Description"
$hh.Range("C1").Value = "Display Value (override)"
$hh.Range("D1").Value = "Last Updated"
$hh.Range("E1").Value = "Note"

# Clone the header formatting (bold font, centered/top alignment, thin
# border) from an existing header cell so we land on the same shared
# style rather than fabricating a new (slightly different) one.
$srcHeader = $wb.Worksheets.Item("Financial_Summary").Range("A1")
$srcHeader.Copy()
$hh.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats

# Data row 2 - only "Last Updated" is populated for now
$hh.Range("D2").Value = "2025-08-10 14:20 UTC"

[void]$hh.Range("A1").Select()
